# "Generate Report for Handoff" - refresh the localization-status report:
#  - status moves from "In Translation" to "Ready for handoff"
#  - the handoff timestamps are bumped to the new generation time
#  - the Status/Datetime columns are widened to fit the new text

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-07-08 09:39:48"
$wsOverview.Columns.Item(5).ColumnWidth = 16.41
$wsOverview.Columns.Item(6).ColumnWidth = 16.41

# --- zh-cn sheet ------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("G2").Value = "2016-07-08 09:39:39"
$wsZh.Columns.Item(3).ColumnWidth = 16.41

# --- de-de sheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("G2").Value = "2016-07-08 09:39:48"
$wsDe.Columns.Item(3).ColumnWidth = 16.41
